$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.150978565216064
$ws.Range("B1").Value = 2.969846725463867
$ws.Range("C1").Value = 3.74387526512146
$ws.Range("D1").Value = 3.467932939529419
$ws.Range("E1").Value = 1.198373198509216
